# Fill in the final (empty) list paragraph with "ExecuteSqlCommand:" and
# append a new list paragraph explaining ExecuteSqlCommand, per the diff:
#   - "ExecuteSqlCommand" (flagged as a spell-check error) + ":" inside the
#     existing empty paragraph (keeps its own ilvl=1 numbering/pPr).
#   - a brand-new ilvl=2 paragraph: "Executa o comando " + bold "DDL/DML" +
#     " fornecido no banco de dados."

$d = $word.ActiveDocument

$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rPrBold = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# --- Step 1: add "ExecuteSqlCommand:" runs into the existing empty last paragraph ---
$lastPara = $d.Paragraphs.Last
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlPart1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
  + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
  + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' `
  + '<w:p>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r>' + $rPr + '<w:t>ExecuteSqlCommand</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r>' + $rPr + '<w:t>:</w:t></w:r>' `
  + '</w:p>' `
  + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xmlPart1)

# --- Step 2: append a brand-new paragraph (ilvl=2) after it, before the sectPr ---
$pPr2 = '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>' + $rPr + '</w:pPr>'

$xmlPart2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
  + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
  + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' `
  + '<w:p>' + $pPr2 `
  + '<w:r>' + $rPr + '<w:t xml:space="preserve">Executa o comando </w:t></w:r>' `
  + '<w:r>' + $rPrBold + '<w:t>DDL/DML</w:t></w:r>' `
  + '<w:r>' + $rPr + '<w:t xml:space="preserve"> fornecido no banco de dados.</w:t></w:r>' `
  + '</w:p>' `
  + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$docEnd = $d.Range($d.Content.End, $d.Content.End)
$docEnd.InsertXML($xmlPart2)
